$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: "374 - Guess Number Higher or Lower"
$ws.Range("A6").Value = " Guess Number Higher or Lower"

# Hyperlink cell (Local Path column) - add hyperlink then apply the
# existing "Hyperlink" cell style so it matches the other rows exactly.
$ws.Hyperlinks.Add($ws.Range("G6"), "374 - Guess Number Higher or Lower", "", "", "374 - Guess Number Higher or Lower")
$ws.Range("G6").Style = "Hyperlink"

$ws.Range("B6").Value = "Interval"
$ws.Range("C6").Value = "No"
$ws.Range("D6").Value = "No"
$ws.Range("E6").Value = "Easy"
$ws.Range("F6").Value = "Easy"

# Extend conditional formatting range to include the new row.
$ws.Range("D2:G5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:G6"))

# Extend the data validation ranges to include the new row, rebuilding
# each rule (in its original order) over the widened range.
$ws.Range("E2:F5").Validation.Delete()
$ws.Range("C2:C5").Validation.Delete()
$ws.Range("B2:B5").Validation.Delete()
$ws.Range("D2:D5").Validation.Delete()

$dv1 = $ws.Range("E2:F6").Validation
$dv1.Add(3, 1, 1, '"Easy, Medium, Hard"')

$dv2 = $ws.Range("C2:C6").Validation
$dv2.Add(3, 1, 1, '"Yes, No"')
$dv2.IgnoreBlank = $false

$dv3 = $ws.Range("B2:B6").Validation
$dv3.Add(3, 1, 1, '"Array, Binary, Dynamic Programming, Graph, Interval, Linked List, Matrix, String, Tree, Heap"')

$dv4 = $ws.Range("D2:D6").Validation
$dv4.Add(3, 1, 1, '"Yes, No"')

# Move the active selection, matching the saved view state.
$ws.Range("R8").Select()
